$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1120.9615
$ws.Range("J17").Value = 1130.2
$ws.Range("L17").Value = 3390.6
$ws.Range("N17").Value = -3726.6
$ws.Range("H64").Value = 3211.9487
$ws.Range("I64").Value = 3003.5293
$ws.Range("J64").Value = 3373
$ws.Range("K64").Value = 3003.5293
$ws.Range("L64").Value = 3373
$ws.Range("M64").Value = -2755.5293
$ws.Range("N64").Value = -3869
$ws.Range("H67").Value = 3211.9487
$ws.Range("I67").Value = 3003.5293
$ws.Range("J67").Value = 3373
$ws.Range("K67").Value = 3003.5293
$ws.Range("L67").Value = 3373
$ws.Range("M67").Value = -2145.5293
$ws.Range("N67").Value = -5089
$ws.Range("H74").Value = 3863
$ws.Range("I74").Value = 3823.8572
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3823.8572
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2887.8572
$ws.Range("N74").Value = -5872
$ws.Range("H76").Value = 3013.5386
$ws.Range("I76").Value = 2839.6
$ws.Range("J76").Value = 3593.3333
$ws.Range("K76").Value = 2839.6
$ws.Range("L76").Value = 3593.3333
$ws.Range("M76").Value = -2524.6
$ws.Range("N76").Value = -4223.3333
$ws.Range("H77").Value = 3863
$ws.Range("I77").Value = 3823.8572
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 19119.286
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -14439.286
$ws.Range("N77").Value = -29360
$ws.Range("H79").Value = 3013.5386
$ws.Range("I79").Value = 2839.6
$ws.Range("J79").Value = 3593.3333
$ws.Range("K79").Value = 2839.6
$ws.Range("L79").Value = 3593.3333
$ws.Range("M79").Value = -1747.6
$ws.Range("N79").Value = -5777.3333
$ws.Range("H103").Value = 553.7222
$ws.Range("I103").Value = 611.4
$ws.Range("J103").Value = 481.625
$ws.Range("K103").Value = 1834.2
$ws.Range("L103").Value = 1444.875
$ws.Range("M103").Value = -1248.2
$ws.Range("N103").Value = -2616.875
$ws.Range("H137").Value = 2235.7144
$ws.Range("I137").Value = 3810
$ws.Range("J137").Value = 1217.0588
$ws.Range("K137").Value = 11430
$ws.Range("L137").Value = 3651.1764
$ws.Range("M137").Value = -8880
$ws.Range("N137").Value = -8751.1764

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1757393.6
$ws.Range("I32").Value = 13112.108
$ws.Range("J32").Value = 17891998
$ws.Range("K32").Value = 13112.108
$ws.Range("L32").Value = 17891998
$ws.Range("M32").Value = -12825.108
$ws.Range("N32").Value = -17892572
$ws.Range("H63").Value = 280311.66
$ws.Range("I63").Value = 335431
$ws.Range("J63").Value = 4715
$ws.Range("K63").Value = 335431
$ws.Range("L63").Value = 4715
$ws.Range("M63").Value = -334745
$ws.Range("N63").Value = -6087
$ws.Range("H66").Value = 280311.66
$ws.Range("I66").Value = 335431
$ws.Range("J66").Value = 4715
$ws.Range("K66").Value = 1677155
$ws.Range("L66").Value = 23575
$ws.Range("M66").Value = -1673723
$ws.Range("N66").Value = -30439
$ws.Range("H88").Value = 1579
$ws.Range("I88").Value = 1358
$ws.Range("J88").Value = 1800
$ws.Range("K88").Value = 1358
$ws.Range("L88").Value = 1800
$ws.Range("M88").Value = -952
$ws.Range("N88").Value = -2612
$ws.Range("H91").Value = 1579
$ws.Range("I91").Value = 1358
$ws.Range("J91").Value = 1800
$ws.Range("K91").Value = 1358
$ws.Range("L91").Value = 1800
$ws.Range("M91").Value = 46
$ws.Range("N91").Value = -4608

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1492.1177
$ws.Range("I86").Value = 1429.9166
$ws.Range("J86").Value = 1641.4
$ws.Range("K86").Value = 1429.9166
$ws.Range("L86").Value = 1641.4
$ws.Range("M86").Value = -306.9166
$ws.Range("N86").Value = -3887.4
$ws.Range("H89").Value = 1492.1177
$ws.Range("I89").Value = 1429.9166
$ws.Range("J89").Value = 1641.4
$ws.Range("K89").Value = 7149.583000000001
$ws.Range("L89").Value = 8207
$ws.Range("M89").Value = -1533.583000000001
$ws.Range("N89").Value = -19439
$ws.Range("H105").Value = 2459.6428
$ws.Range("I105").Value = 2265.238
$ws.Range("J105").Value = 3042.8572
$ws.Range("K105").Value = 2265.238
$ws.Range("L105").Value = 3042.8572
$ws.Range("M105").Value = -518.2379999999998
$ws.Range("N105").Value = -6536.8572
$ws.Range("H130").Value = 35000
$ws.Range("J130").Value = 35000
$ws.Range("L130").Value = 35000
$ws.Range("N130").Value = -45040

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3062.6936
$ws.Range("I31").Value = 1714.65
$ws.Range("J31").Value = 5513.6816
$ws.Range("K31").Value = 1714.65
$ws.Range("L31").Value = 5513.6816
$ws.Range("M31").Value = -1419.65
$ws.Range("N31").Value = -6103.6816
$ws.Range("H34").Value = 3062.6936
$ws.Range("I34").Value = 1714.65
$ws.Range("J34").Value = 5513.6816
$ws.Range("K34").Value = 1714.65
$ws.Range("L34").Value = 5513.6816
$ws.Range("M34").Value = -1512.65
$ws.Range("N34").Value = -5917.6816
$ws.Range("H62").Value = 11336.091
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 34899
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 34899
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -36147
$ws.Range("H65").Value = 11336.091
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 34899
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 174495
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -180735
$ws.Range("H132").Value = 2046.0817
$ws.Range("I132").Value = 1776.7368
$ws.Range("J132").Value = 2216.6667
$ws.Range("K132").Value = 5330.2104
$ws.Range("L132").Value = 6650.000100000001
$ws.Range("M132").Value = -2800.2104
$ws.Range("N132").Value = -11710.0001
$ws.Range("H134").Value = 1826.9302
$ws.Range("I134").Value = 1068.8334
$ws.Range("J134").Value = 2372.76
$ws.Range("K134").Value = 3206.5002
$ws.Range("L134").Value = 7118.280000000001
$ws.Range("M134").Value = -671.5001999999999
$ws.Range("N134").Value = -12188.28
$ws.Range("H135").Value = 46100
$ws.Range("J135").Value = 46100
$ws.Range("L135").Value = 46100
$ws.Range("N135").Value = -56240

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2699.6667
$ws.Range("J80").Value = 2649.611
$ws.Range("L80").Value = 7948.833
$ws.Range("N80").Value = -9820.832999999999
$ws.Range("H83").Value = 2699.6667
$ws.Range("J83").Value = 2649.611
$ws.Range("L83").Value = 23846.499
$ws.Range("N83").Value = -33206.499
$ws.Range("H113").Value = 22223106
$ws.Range("I113").Value = 656.6
$ws.Range("J113").Value = 33334330
$ws.Range("K113").Value = 1969.8
$ws.Range("L113").Value = 100002990
$ws.Range("M113").Value = 200.1999999999998
$ws.Range("N113").Value = -100007330
$ws.Range("H117").Value = 35717280
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 35717280
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 107151840
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -107158724
$ws.Range("H129").Value = 1763.7059
$ws.Range("I129").Value = 822.8570999999999
$ws.Range("K129").Value = 2468.5713
$ws.Range("M129").Value = 2531.4287

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6201.3335
$ws.Range("I70").Value = 7314.2856
$ws.Range("J70").Value = 5493.091
$ws.Range("K70").Value = 7314.2856
$ws.Range("L70").Value = 5493.091
$ws.Range("M70").Value = -7044.2856
$ws.Range("N70").Value = -6033.091
$ws.Range("H73").Value = 6201.3335
$ws.Range("I73").Value = 7314.2856
$ws.Range("J73").Value = 5493.091
$ws.Range("K73").Value = 7314.2856
$ws.Range("L73").Value = 5493.091
$ws.Range("M73").Value = -6378.2856
$ws.Range("N73").Value = -7365.091
$ws.Range("H80").Value = 2860.2666
$ws.Range("I80").Value = 2650
$ws.Range("J80").Value = 2936.7273
$ws.Range("K80").Value = 2650
$ws.Range("L80").Value = 2936.7273
$ws.Range("M80").Value = -1652
$ws.Range("N80").Value = -4932.7273
$ws.Range("H83").Value = 2860.2666
$ws.Range("I83").Value = 2650
$ws.Range("J83").Value = 2936.7273
$ws.Range("K83").Value = 13250
$ws.Range("L83").Value = 14683.6365
$ws.Range("M83").Value = -8258
$ws.Range("N83").Value = -24667.6365

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 31253308
$ws.Range("I40").Value = 37040004
$ws.Range("J40").Value = 5157
$ws.Range("K40").Value = 37040004
$ws.Range("L40").Value = 5157
$ws.Range("M40").Value = -37039868
$ws.Range("N40").Value = -5429

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1380.9595
$ws.Range("I132").Value = 1037.7
$ws.Range("J132").Value = 2096.0833
$ws.Range("K132").Value = 3113.1
$ws.Range("L132").Value = 6288.249899999999
$ws.Range("M132").Value = -583.1000000000004
$ws.Range("N132").Value = -11348.2499
